$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C13").Value = " Browsing is just relaxing!"
$ws.Range("A13").Value = "SCRIPT/P01P04A/um1203.ssb"
$ws.Range("D13").Value = " Ходьба расслабляет!"
$ws.Range("E13").Value = " Öïäûáà ñàòòìàáìÿåó!"
$ws.Range("B13").Value = 239
